$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1 - copy E1's header style (bold, centered, bordered) and set text value
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("F1").Value = "time_taken"

$timestamps = @(
    "2021-10-05 10:52:34.528361",
    "2021-10-05 10:52:34.528370",
    "2021-10-05 10:52:34.528373",
    "2021-10-05 10:52:34.528376",
    "2021-10-05 10:52:34.528379",
    "2021-10-05 10:52:34.528381",
    "2021-10-05 10:52:34.528384",
    "2021-10-05 10:52:34.528386",
    "2021-10-05 10:52:34.528389",
    "2021-10-05 10:52:34.528392",
    "2021-10-05 10:52:34.528394",
    "2021-10-05 10:52:34.528396",
    "2021-10-05 10:52:34.528399",
    "2021-10-05 10:52:34.528401",
    "2021-10-05 10:52:34.528404",
    "2021-10-05 10:52:34.528406",
    "2021-10-05 10:52:34.528409"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
